$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" values (column G) replacing the old "Strike#" values,
# recomputed from the regenerated save_data (std/mean recalculated, s_vals written).
$newK = @{
    2  = 0
    3  = 0
    4  = 2
    5  = 0
    6  = 2
    7  = 0
    8  = 1
    9  = 0
    10 = 3
    11 = 1
    12 = 2
    13 = 1
    14 = 1
    16 = 2
    17 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
